$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look like plain numbers,
# so Excel keeps them as text (matching the original text-formatted price column).
$textCells = @("D5", "D6", "D8", "D10", "D11", "D17", "D20", "D21", "D26", "D27", "D38", "D43", "D44", "D45", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values from the crypto price refresh.
$ws.Range("D2").Value = "25.768.53"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "1.634.65"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "215.68"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").Value = "0.500"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").Value = "0.256"
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("E9").Value = "  -1.29%  "
$ws.Range("D10").Value = "19.59"
$ws.Range("E10").Value = "  -1.37%  "
$ws.Range("D11").Value = "0.0792"
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").Value = "1.860.35"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").Value = "1.643.32"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("D16").Value = "0.0₃0763"
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").Value = "63.20"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "25.794.31"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").Value = "4.46"
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("D21").Value = "192.51"
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("E23").Value = "  +2.45%  "
$ws.Range("E24").Value = "  +3.64%  "
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").Value = "142.13"
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("D27").Value = "0.124"
$ws.Range("E27").Value = "  +2.01%  "
$ws.Range("E28").Value = "  +0.92%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("E33").Value = "  -0.71%  "
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").Value = "1.132.18"
$ws.Range("E37").Value = "  +1.78%  "
$ws.Range("D38").Value = "0.544"
$ws.Range("E38").Value = "  -1.37%  "
$ws.Range("E40").Value = "  -1.33%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  -0.39%  "
$ws.Range("D43").Value = "5.55"
$ws.Range("E43").Value = "  -0.47%  "
$ws.Range("D44").Value = "100.75"
$ws.Range("E44").Value = "  +1.15%  "
$ws.Range("D45").Value = "0.800"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("D46").Value = "1.769.63"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("E47").Value = "  +0.66%  "
$ws.Range("D48").Value = "55.42"
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0505"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "0.417"
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("E51").Value = "  +3.78%  "
